$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

# Fix header P1: remove trailing space from "Event "
$ws.Range("P1").Value = "Event"

# Add new header Q1: "Correction " (with trailing space), matching style/format of P1
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("Q1").Value = "Correction "

# Fill column P (rows 2-12) with "nan" where currently blank
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 16).Value = "nan"
}

# Column Q (rows 2-12) stays blank, but the cells still need to exist in
# the sheet (present-but-empty), matching the source diff. Touching a
# formatting property with its already-default value materializes the
# cell record without introducing a new style or any value.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 17).Font.Bold = $false
}
